$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new header row (row 4): Name / Form / To / Paid ---
$ws.Range("C4").Value = "Name"
$ws.Range("D4").Value = "Form"
$ws.Range("E4").Value = "To"
$ws.Range("F4").Value = "Paid"

# Center the new row, same alignment style already used elsewhere in the sheet.
$ws.Range("C4:F4").HorizontalAlignment = -4108
Write-Output "row 4 written"

# --- Re-stamp row 3 (the merged title row) with its own distinct (but
# visually identical - centered, no border) style so it no longer shares
# the exact same style index as the newly written row 4. ---
$ws.Range("C3:F3").Locked = $true
$ws.Range("C3:F3").HorizontalAlignment = -4108
Write-Output "row 3 restyled"

# --- Move the active selection from E9 to E8 ---
$null = $ws.Range("E8").Select()
Write-Output "selection updated"
